# "Add files via upload" — the sheet now carries a single text value.
# Reproduce: Sheet1!A1 = "cccccc", leaving the selection on K4 afterward
# (matches the <selection activeCell="K4" sqref="K4"/> recorded in the
# saved sheetView of the edited workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "cccccc"

$ws.Range("K4").Select() | Out-Null
